# Mise à jour de l'application
# Adds 12 new training-log rows (142-153) for the 2025-08-14 session
# to the "Wellness" worksheet, extending the existing table/formula range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date serial (2025-08-14)
$newDate = 45883

# Data for the new rows:
# RowOffset, PlayerName, Volume, Intensite, Fatigue, Douleur, Localisation, Plaisir
$rows = @(
    @{ Player = "Amir Etien";       Volume = 60; Intensite = 5; Fatigue = 6; Douleur = 0; Loc = "";          Plaisir = 6  },
    @{ Player = "Maé Clavel";       Volume = 60; Intensite = 3; Fatigue = 4; Douleur = 3; Loc = "Cheville";  Plaisir = 6  },
    @{ Player = "Ilyes Boughanmi";  Volume = 60; Intensite = 6; Fatigue = 5; Douleur = 6; Loc = "Ampoule";   Plaisir = 10 },
    @{ Player = "Omar Benyounes";   Volume = 60; Intensite = 5; Fatigue = 5; Douleur = 0; Loc = "";          Plaisir = 7  },
    @{ Player = "Naim Ighbane";     Volume = 60; Intensite = 6; Fatigue = 6; Douleur = 0; Loc = "";          Plaisir = 2  },
    @{ Player = "Karim Belmahi";    Volume = 60; Intensite = 6; Fatigue = 7; Douleur = 0; Loc = "";          Plaisir = 10 },
    @{ Player = "Rayane Chayebi";   Volume = 60; Intensite = 6; Fatigue = 6; Douleur = 4; Loc = "Adducteurs";Plaisir = 7  },
    @{ Player = "Romain Thunet";    Volume = 60; Intensite = 6; Fatigue = 4; Douleur = 4; Loc = "Genou";     Plaisir = 2  },
    @{ Player = "Ilan Ihaddadene";  Volume = 60; Intensite = 6; Fatigue = 7; Douleur = 0; Loc = "";          Plaisir = 8  },
    @{ Player = "Naim Dhib";        Volume = 60; Intensite = 6; Fatigue = 5; Douleur = 0; Loc = "";          Plaisir = 0  },
    @{ Player = "Mattheo Haon";     Volume = 60; Intensite = 5; Fatigue = 7; Douleur = 0; Loc = "";          Plaisir = 5  },
    @{ Player = "Yanis Berrached";  Volume = 60; Intensite = 6; Fatigue = 7; Douleur = 3; Loc = "Courbature";Plaisir = 7  }
)

$startRow = 142
$lastOldRow = 141
$endRow = $startRow + $rows.Count - 1

# Extend the sheet by copying formatting from the last existing row down
# to all the new rows in one shot.
$ws.Range("A$($lastOldRow):I$($lastOldRow)").Copy()
$ws.Range("A$($startRow):I$($endRow)").PasteSpecial(-4122)

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 2).Value = $row.Player
    $ws.Cells.Item($r, 3).Value = $row.Volume
    $ws.Cells.Item($r, 4).Value = $row.Intensite
    $ws.Cells.Item($r, 5).Value = $row.Fatigue
    $ws.Cells.Item($r, 6).Value = $row.Douleur

    if ($row.Loc -ne "") {
        # Match the formatting used on other rows that have a "Localisation" value
        $ws.Range("G134").Copy()
        $ws.Cells.Item($r, 7).PasteSpecial(-4122)
        $ws.Cells.Item($r, 7).Value = $row.Loc
    }

    $ws.Cells.Item($r, 8).Value = $row.Plaisir

    $r++
}

# Charge column (I) = Volume * Intensite, carried down as a formula like the rest of the table
$ws.Range("I$($startRow):I$($endRow)").Formula = "=C$($startRow)*D$($startRow)"

# Update the view to match where the workbook was scrolled to after the edit
$ws.Range("K144").Select()
